# The original workbook has a single empty sheet named "Sheet1".
# This edit turns it into a 3-sheet data-driven-test workbook:
#   test_suite, AddCustomerTest, OpenAccountTest   (tab order, left to right)
#
# Step 1 - do all the *structural* sheet work first (rename the
# original sheet, add the two new ones, then reorder tabs) before any
# cell is touched. Writing data only after the final layout is in
# place keeps each worksheet's on-disk identity stable and lets the
# physical file numbering line up with creation order.

$wb = $excel.ActiveWorkbook

# Original (only) sheet becomes "AddCustomerTest".
$wsAddCustomer = $wb.Worksheets.Item(1)
$wsAddCustomer.Name = "AddCustomerTest"

# New sheet "OpenAccountTest" right after it.
$wsOpenAccount = $wb.Worksheets.Add($null, $wsAddCustomer)
$wsOpenAccount.Name = "OpenAccountTest"

# New sheet "test_suite" right after that.
$wsTestSuite = $wb.Worksheets.Add($null, $wsOpenAccount)
$wsTestSuite.Name = "test_suite"

# Reorder tabs so test_suite leads: test_suite, AddCustomerTest, OpenAccountTest.
$wsTestSuite.Move($wb.Worksheets.Item(1))

# Step 2 - populate each sheet with its data / view state.

# ---------------- AddCustomerTest ----------------
$wsAddCustomer.Cells.Item(1, 1).Value = "TCID"
$wsAddCustomer.Cells.Item(1, 2).Value = "Runmode"
$wsAddCustomer.Cells.Item(2, 1).Value = "BankManagerLoginTest"
$wsAddCustomer.Cells.Item(2, 2).Value = "Y"
$wsAddCustomer.Cells.Item(3, 1).Value = "AddCustomerTest"
$wsAddCustomer.Cells.Item(3, 2).Value = "Y"
$wsAddCustomer.Cells.Item(4, 1).Value = "OpenAccountTest"
$wsAddCustomer.Cells.Item(4, 2).Value = "y"
# Column A was widened (auto-fit) to show the longest TCID name in full.
$wsAddCustomer.Columns.Item(1).ColumnWidth = 19.666666666666668
$wsAddCustomer.Range("B5").Select()

# ---------------- OpenAccountTest ----------------
$wsOpenAccount.Cells.Item(1, 1).Value = "firstName"
$wsOpenAccount.Cells.Item(1, 2).Value = "lastName"
$wsOpenAccount.Cells.Item(1, 3).Value = "postCode"
$wsOpenAccount.Cells.Item(1, 4).Value = "alertText"
$wsOpenAccount.Cells.Item(1, 5).Value = "runmode"

$wsOpenAccount.Cells.Item(2, 1).Value = "Raman"
$wsOpenAccount.Cells.Item(2, 2).Value = "Arora"
$wsOpenAccount.Cells.Item(2, 3).Value = 1212
$wsOpenAccount.Cells.Item(2, 4).Value = "Customer Added successfully"
$wsOpenAccount.Cells.Item(2, 5).Value = "y"

$wsOpenAccount.Cells.Item(3, 1).Value = "Rahul"
$wsOpenAccount.Cells.Item(3, 2).Value = "Arora"
$wsOpenAccount.Cells.Item(3, 3).Value = 1212
$wsOpenAccount.Cells.Item(3, 4).Value = "Customer Added successfully"
$wsOpenAccount.Cells.Item(3, 5).Value = "N"

$wsOpenAccount.Cells.Item(4, 1).Value = "Ishita"
$wsOpenAccount.Cells.Item(4, 2).Value = "Arora"
$wsOpenAccount.Cells.Item(4, 3).Value = 1212
$wsOpenAccount.Cells.Item(4, 4).Value = "Customer Added successfully"
$wsOpenAccount.Cells.Item(4, 5).Value = "y"

$wsOpenAccount.Cells.Item(5, 1).Value = "Rohit"
$wsOpenAccount.Cells.Item(5, 2).Value = "Sehgal"
$wsOpenAccount.Cells.Item(5, 3).Value = 1212
$wsOpenAccount.Cells.Item(5, 4).Value = "Customer Added successfully"
$wsOpenAccount.Cells.Item(5, 5).Value = "y"

# Column D was widened (auto-fit) to show the full alert text.
$wsOpenAccount.Columns.Item(4).ColumnWidth = 24.3
$wsOpenAccount.Range("E5").Select()

# ---------------- test_suite ----------------
$wsTestSuite.Cells.Item(1, 1).Value = "customer"
$wsTestSuite.Cells.Item(1, 2).Value = "currency"
$wsTestSuite.Cells.Item(2, 1).Value = "Raman Arora"
$wsTestSuite.Cells.Item(2, 2).Value = "Rupee"
# Column A was widened (auto-fit) to show "customer" / "Raman Arora" in full.
$wsTestSuite.Columns.Item(1).ColumnWidth = 11.166666666666668
$wsTestSuite.Columns.Item(2).Select()

# OpenAccountTest was the sheet on screen when the workbook was saved.
$wsOpenAccount.Select()
